$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (values re-scraped by the GitHub Actions job).
$updates = [ordered]@{
    "D2" = "297.55"
    "E2" = "2.90%"
    "E3" = "2.34%"
    "D4" = "5.011"
    "E4" = "-0.57%"
    "D5" = "0.07534"
    "E5" = "3.40%"
    "D6" = "4.366"
    "E6" = "1.91%"
    "D7" = "1.573"
    "E7" = "3.05%"
    "D8" = "0.9307"
    "E8" = "1.35%"
    "D10" = "0.1199"
    "E10" = "2.05%"
    "D11" = "0.1822"
    "E11" = "5.93%"
    "D12" = "0.08824"
    "E12" = "1.52%"
    "D13" = "0.04079"
    "E13" = "-2.33%"
    "D14" = "0.1055"
    "E14" = "0.17%"
    "D15" = "0.001279"
    "E15" = "0.36%"
    "D16" = "0.005917"
    "E16" = "1.51%"
    "E17" = "-1.40%"
    "D18" = "0.3335"
    "E18" = "0.50%"
    "D19" = "7.921"
    "E19" = "0.42%"
    "E20" = "4.75%"
    "D21" = "0.2995"
    "E21" = "3.75%"
    "D22" = "0.04055"
    "E22" = "5.07%"
    "E23" = "-0.64%"
    "D24" = "0.003895"
    "E24" = "1.09%"
    "D25" = "0.0001229"
    "E25" = "-4.07%"
    "D38" = "0.02419"
    "E38" = "4.59%"
    "E39" = "5.25%"
    "D40" = "0.006071"
    "E40" = "-7.42%"
    "D41" = "0.007784"
    "E41" = "1.30%"
    "D42" = "0.1333"
    "E42" = "4.66%"
    "D43" = "0.007356"
    "E43" = "0.04%"
    "D44" = "0.007832"
    "E44" = "10.86%"
    "D45" = "0.2977"
    "E45" = "-4.78%"
    "D46" = "0.00006302"
    "E46" = "-2.30%"
    "E47" = "-0.34%"
    "D48" = "0.04326"
    "E48" = "406.61%"
    "D49" = "0.004194"
    "E49" = "-0.18%"
    "D50" = "0.00002097"
    "E50" = "-0.34%"
    "D51" = "0.0001997"
    "E51" = "-0.34%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (and the literal "%" suffix
    # on the Volume(1h) column) round-trip as text instead of Excel auto-detecting
    # Number/Percentage types for them.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
